$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 8080  # was 8062
$ws.Range("F5").Value = 943  # was 942
$ws.Range("F8").Value = 622  # was 617
$ws.Range("F9").Value = 103  # was 101
$ws.Range("F12").Value = 873  # was 872
$ws.Range("F13").Value = 3317  # was 3306
$ws.Range("F14").Value = 220  # was 215
$ws.Range("F15").Value = 108  # was 107
$ws.Range("F16").Value = 752  # was 748
$ws.Range("F17").Value = 760  # was 758
$ws.Range("F19").Value = 465  # was 464
$ws.Range("F21").Value = 279  # was 275
$ws.Range("F22").Value = 487  # was 236
$ws.Range("F23").Value = 352  # was 348
$ws.Range("F25").Value = 134  # was 132
$ws.Range("F26").Value = 128  # was 125
$ws.Range("F27").Value = 294  # was 288
$ws.Range("F28").Value = 33  # was 32
$ws.Range("F32").Value = 579  # was 569
$ws.Range("F33").Value = 28  # was 27
$ws.Range("F34").Value = 37  # was 36
$ws.Range("F35").Value = 20  # was 17
$ws.Range("F38").Value = 108  # was 107

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 21  # was 20

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 8080  # was 8062
$ws.Range("F7").Value = 943  # was 942
$ws.Range("F10").Value = 622  # was 617
$ws.Range("F11").Value = 103  # was 101
$ws.Range("F14").Value = 873  # was 872
$ws.Range("F15").Value = 21  # was 20
$ws.Range("F16").Value = 3317  # was 3307
$ws.Range("F17").Value = 220  # was 215
$ws.Range("F18").Value = 108  # was 107
$ws.Range("F20").Value = 752  # was 748
$ws.Range("F21").Value = 760  # was 758
$ws.Range("F24").Value = 465  # was 464
$ws.Range("F26").Value = 279  # was 275
$ws.Range("F27").Value = 489  # was 236
$ws.Range("F28").Value = 352  # was 349
$ws.Range("F30").Value = 134  # was 133
$ws.Range("F31").Value = 128  # was 125
$ws.Range("F32").Value = 294  # was 288
$ws.Range("F33").Value = 33  # was 32
$ws.Range("F37").Value = 579  # was 569
$ws.Range("F38").Value = 28  # was 27
$ws.Range("F39").Value = 37  # was 36
$ws.Range("F40").Value = 20  # was 17
$ws.Range("F43").Value = 108  # was 107
